# Center the text in the four "Rounded Rectangle" shapes on slide 23
# (Asynchrony, Migratability, Introspection, Adaptivity) to match the
# alignment already used by their neighboring shapes (Overdecomposition,
# Adaptive / Runtime System).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(23)

$shapeNames = @("Rounded Rectangle 7", "Rounded Rectangle 9", "Rounded Rectangle 11", "Rounded Rectangle 12")

foreach ($name in $shapeNames) {
    $sh = $s.Shapes.Item($name)
    $sh.TextFrame.TextRange.ParagraphFormat.Alignment = 2  # ppAlignCenter
}
